# Version 3 of the generated data: refreshed tripleUuid identifiers for the
# Chemicals & Drugs / Anatomy / Disorders tiers, and recomputed pathWeight
# scores for the two rows that were previously missing.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Chemicals & Drugs")
# tier01TripleInformation/0/tripleUuid (column H): force text format while
# writing so the numeric-looking IDs are kept as strings, then clear the
# temporary format back off again
$uuidRange = $ws.Range("H2:H15")
$uuidRange.NumberFormat = "@"
$ws.Range("H2").Value = "76350720"
$ws.Range("H3").Value = "78841912"
$ws.Range("H4").Value = "116859320"
$ws.Range("H5").Value = "74231788"
$ws.Range("H6").Value = "58409705"
$ws.Range("H7").Value = "56773882"
$ws.Range("H8").Value = "75949843"
$ws.Range("H9").Value = "114043056"
$ws.Range("H10").Value = "61382531"
$ws.Range("H11").Value = "94691585"
$ws.Range("H12").Value = "58031591"
$ws.Range("H13").Value = "53339440"
$ws.Range("H14").Value = "125280756"
$ws.Range("H15").Value = "114040142"
$uuidRange.ClearFormats()
# pathWeight (column A)
$ws.Range("A4").Value = 5.313835620880127
$ws.Range("A8").Value = 5.013763904571533
$ws.Range("A12").Value = 4.692009449005127

$ws = $wb.Worksheets.Item("Anatomy")
# tier01TripleInformation/0/tripleUuid (column H): force text format while
# writing so the numeric-looking IDs are kept as strings, then clear the
# temporary format back off again
$uuidRange = $ws.Range("H2:H16")
$uuidRange.NumberFormat = "@"
$ws.Range("H2").Value = "58031670"
$ws.Range("H3").Value = "59479264"
$ws.Range("H4").Value = "70409366"
$ws.Range("H5").Value = "57099226"
$ws.Range("H6").Value = "59517342"
$ws.Range("H7").Value = "121751908"
$ws.Range("H8").Value = "78783751"
$ws.Range("H9").Value = "87816211"
$ws.Range("H10").Value = "53682816"
$ws.Range("H11").Value = "102897923"
$ws.Range("H12").Value = "108166559"
$ws.Range("H13").Value = "121941152"
$ws.Range("H14").Value = "122848718"
$ws.Range("H15").Value = "130743106"
$ws.Range("H16").Value = "122851169"
$uuidRange.ClearFormats()

$ws = $wb.Worksheets.Item("Disorders")
# tier01TripleInformation/0/tripleUuid (column H): force text format while
# writing so the numeric-looking IDs are kept as strings, then clear the
# temporary format back off again
$uuidRange = $ws.Range("H2:H20")
$uuidRange.NumberFormat = "@"
$ws.Range("H2").Value = "113483617"
$ws.Range("H3").Value = "127281809"
$ws.Range("H4").Value = "57273052"
$ws.Range("H5").Value = "68296895"
$ws.Range("H6").Value = "84843954"
$ws.Range("H7").Value = "55291627"
$ws.Range("H8").Value = "137260534"
$ws.Range("H9").Value = "132797091"
$ws.Range("H10").Value = "66026027"
$ws.Range("H11").Value = "103628962"
$ws.Range("H12").Value = "131050434"
$ws.Range("H13").Value = "68735561"
$ws.Range("H14").Value = "60580724"
$ws.Range("H15").Value = "110982258"
$ws.Range("H16").Value = "88555167"
$ws.Range("H17").Value = "116859123"
$ws.Range("H18").Value = "63359632"
$ws.Range("H19").Value = "138413405"
$ws.Range("H20").Value = "71155385"
$uuidRange.ClearFormats()
# pathWeight (column A)
$ws.Range("A2").Value = 5.908326625823975
$ws.Range("A5").Value = 5.6477155685424805
$ws.Range("A6").Value = 5.51275110244751
$ws.Range("A7").Value = 5.342600345611572
